# Applies:
#  1. Date placeholder text "7/16/2019" -> "8/23/2019" on the Slide Master
#     and every Slide Layout (the "datetimeFigureOut" field placeholder).
#  2. Slide 6 "TextBox 6": "Assumes catch in 100% mature." -> "Assumes catch is 100% mature."

$p = $ppt.ActivePresentation

function Get-DatePlaceholderShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            return $sh
        }
    }
    return $null
}

function Update-DateShapeText($sh, $oldText, $newText) {
    if ($sh -eq $null) {
        return
    }
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -eq $oldText) {
        $chars = $tr.Characters(1, $oldText.Length)
        $chars.Text = $newText
    }
}

$oldDate = "7/16/2019"
$newDate = "8/23/2019"

$master = $p.SlideMaster

# Slide Master date placeholder.
Update-DateShapeText (Get-DatePlaceholderShape $master) $oldDate $newDate

# Every Slide Layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapeText (Get-DatePlaceholderShape $layout) $oldDate $newDate
}

# Slide 6: fix "catch in 100%" -> "catch is 100%" in "TextBox 6".
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $sh = $slide6.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 6") {
        $tr = $sh.TextFrame.TextRange
        $oldRun = "Assumes catch in 100% mature. Uses "
        $newRun = "Assumes catch is 100% mature. Uses "
        if ($tr.Text.StartsWith($oldRun)) {
            $chars = $tr.Characters(1, $oldRun.Length)
            $chars.Text = $newRun
        }
    }
}
